$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2024-04-19 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-20 Saturday", 2) | Out-Null

# Update the division-problem answer table (positional addressing
# avoids ambiguity from values that collide with other cells' before/after text)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "25÷8=3, 1"
$t.Cell(1, 2).Range.Text = "37÷9=4, 1"
$t.Cell(1, 3).Range.Text = "17÷7=2, 3"
$t.Cell(1, 4).Range.Text = "54÷5=10, 4"
$t.Cell(1, 5).Range.Text = "90÷6=15, 0"
$t.Cell(5, 1).Range.Text = "65÷2=32, 1"
$t.Cell(5, 2).Range.Text = "51÷5=10, 1"
$t.Cell(5, 3).Range.Text = "52÷9=5, 7"
$t.Cell(5, 4).Range.Text = "23÷5=4, 3"
$t.Cell(5, 5).Range.Text = "51÷8=6, 3"
$t.Cell(9, 1).Range.Text = "24÷8=3, 0"
$t.Cell(9, 2).Range.Text = "94÷9=10, 4"
$t.Cell(9, 3).Range.Text = "99÷6=16, 3"
$t.Cell(9, 4).Range.Text = "33÷9=3, 6"
$t.Cell(9, 5).Range.Text = "95÷9=10, 5"
$t.Cell(13, 1).Range.Text = "45÷6=7, 3"
$t.Cell(13, 2).Range.Text = "55÷8=6, 7"
$t.Cell(13, 3).Range.Text = "96÷7=13, 5"
$t.Cell(13, 4).Range.Text = "87÷7=12, 3"
$t.Cell(13, 5).Range.Text = "75÷2=37, 1"
$t.Cell(17, 1).Range.Text = "57÷9=6, 3"
$t.Cell(17, 2).Range.Text = "30÷5=6, 0"
$t.Cell(17, 3).Range.Text = "27÷3=9, 0"
$t.Cell(17, 4).Range.Text = "44÷7=6, 2"
$t.Cell(17, 5).Range.Text = "61÷9=6, 7"
